$wb = $excel.ActiveWorkbook

# ===================== ALC =====================
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H2").Value = 886.3684
$ws.Range("I2").Value = 336.66666
$ws.Range("J2").Value = 1828.7142
$ws.Range("K2").Value = 336.66666
$ws.Range("L2").Value = 1828.7142
$ws.Range("M2").Value = -223.66666
$ws.Range("N2").Value = -2054.7142

$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = $null

$ws.Range("H5").Value = 573.6667
$ws.Range("I5").Value = 566.7273
$ws.Range("K5").Value = 566.7273
$ws.Range("M5").Value = -451.7273

$ws.Range("H17").Value = 2227.8
$ws.Range("J17").Value = 2227.8
$ws.Range("L17").Value = 6683.400000000001
$ws.Range("N17").Value = -7019.400000000001

$ws.Range("H41").Value = 2718
$ws.Range("I41").Value = 2718
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 2718
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -2278
$ws.Range("N41").Value = $null

$ws.Range("H58").Value = 850
$ws.Range("I58").Value = 850
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 2550
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -2400
$ws.Range("N58").Value = $null

$ws.Range("H61").Value = 2671.1428
$ws.Range("I61").Value = 1139.6
$ws.Range("K61").Value = 3418.8
$ws.Range("M61").Value = -3246.8

$ws.Range("H76").Value = 7348.5454
$ws.Range("I76").Value = 3992
$ws.Range("J76").Value = 8094.4443
$ws.Range("K76").Value = 3992
$ws.Range("L76").Value = 8094.4443
$ws.Range("M76").Value = -3677
$ws.Range("N76").Value = -8724.444299999999

$ws.Range("H79").Value = 7348.5454
$ws.Range("I79").Value = 3992
$ws.Range("J79").Value = 8094.4443
$ws.Range("K79").Value = 3992
$ws.Range("L79").Value = 8094.4443
$ws.Range("M79").Value = -2900
$ws.Range("N79").Value = -10278.4443

$ws.Range("H86").Value = 15000
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 15000
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 15000
$ws.Range("M86").Value = $null
$ws.Range("N86").Value = -17246

$ws.Range("H87").Value = 71374.875
$ws.Range("J87").Value = 71374.875
$ws.Range("L87").Value = 71374.875
$ws.Range("N87").Value = -73870.875

$ws.Range("H89").Value = 15000
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 15000
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 75000
$ws.Range("M89").Value = $null
$ws.Range("N89").Value = -86232

$ws.Range("H90").Value = 71374.875
$ws.Range("J90").Value = 71374.875
$ws.Range("L90").Value = 214124.625
$ws.Range("N90").Value = -226604.625

$ws.Range("H98").Value = 1182.0769
$ws.Range("I98").Value = 1182.0769
$ws.Range("K98").Value = 1182.0769
$ws.Range("M98").Value = 315.9231

$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = $null

$ws.Range("H112").Value = 5223
$ws.Range("J112").Value = 5625.5557
$ws.Range("L112").Value = 16876.6671
$ws.Range("N112").Value = -19092.6671

$ws.Range("H113").Value = 4428.1665
$ws.Range("I113").Value = 4253.8
$ws.Range("J113").Value = 5300
$ws.Range("K113").Value = 4253.8
$ws.Range("L113").Value = 5300
$ws.Range("M113").Value = -999.8000000000002
$ws.Range("N113").Value = -11808

$ws.Range("H116").Value = 19142.857
$ws.Range("I116").Value = 19126.25
$ws.Range("K116").Value = 19126.25
$ws.Range("M116").Value = -15684.25

$ws.Range("H121").Value = 5350.6665
$ws.Range("J121").Value = 5350.6665
$ws.Range("L121").Value = 16051.9995
$ws.Range("N121").Value = -19545.9995

$ws.Range("H122").Value = 1182.0769
$ws.Range("I122").Value = 1182.0769
$ws.Range("K122").Value = 3546.2307
$ws.Range("M122").Value = -1096.2307

$ws.Range("H132").Value = 2720.75
$ws.Range("I132").Value = 2674.1853
$ws.Range("J132").Value = 2972.2
$ws.Range("K132").Value = 8022.5559
$ws.Range("L132").Value = 8916.599999999999
$ws.Range("M132").Value = -5492.5559
$ws.Range("N132").Value = -13976.6

$ws.Range("H134").Value = 76458
$ws.Range("J134").Value = 76458
$ws.Range("L134").Value = 76458
$ws.Range("N134").Value = -86598

$ws.Range("H137").Value = 1735.5491
$ws.Range("I137").Value = 1679.4584
$ws.Range("J137").Value = 2633
$ws.Range("K137").Value = 5038.3752
$ws.Range("L137").Value = 7899
$ws.Range("M137").Value = -2488.3752
$ws.Range("N137").Value = -12999


# ===================== ARM =====================
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H2").Value = 1395.4286
$ws.Range("I2").Value = 1445.1538
$ws.Range("J2").Value = 749
$ws.Range("K2").Value = 1445.1538
$ws.Range("L2").Value = 749
$ws.Range("M2").Value = -1332.1538
$ws.Range("N2").Value = -975

$ws.Range("H28").Value = 16749.25
$ws.Range("I28").Value = 16749.25
$ws.Range("K28").Value = 16749.25
$ws.Range("M28").Value = -16557.25

$ws.Range("H32").Value = 8144.803
$ws.Range("I32").Value = 4362.577
$ws.Range("K32").Value = 4362.577
$ws.Range("M32").Value = -4075.577

$ws.Range("H45").Value = 2186.36
$ws.Range("I45").Value = 1172.9375
$ws.Range("J45").Value = 3988
$ws.Range("K45").Value = 1172.9375
$ws.Range("L45").Value = 3988
$ws.Range("M45").Value = -795.9375
$ws.Range("N45").Value = -4742

$ws.Range("H46").Value = 16076
$ws.Range("J46").Value = 16076
$ws.Range("L46").Value = 16076
$ws.Range("N46").Value = -16714

$ws.Range("H74").Value = 31255750
$ws.Range("I74").Value = 41668416
$ws.Range("K74").Value = 41668416
$ws.Range("M74").Value = -41667542

$ws.Range("H77").Value = 31255750
$ws.Range("I77").Value = 41668416
$ws.Range("K77").Value = 208342080
$ws.Range("M77").Value = -208337712

$ws.Range("H99").Value = 16749.25
$ws.Range("I99").Value = 16749.25
$ws.Range("K99").Value = 16749.25
$ws.Range("M99").Value = -13754.25

$ws.Range("H102").Value = 1340.2778
$ws.Range("I102").Value = 1340.2778
$ws.Range("K102").Value = 1340.2778
$ws.Range("M102").Value = 281.7221999999999

$ws.Range("H116").Value = 1395.4286
$ws.Range("I116").Value = 1445.1538
$ws.Range("J116").Value = 749
$ws.Range("K116").Value = 1445.1538
$ws.Range("L116").Value = 749
$ws.Range("M116").Value = 848.8462
$ws.Range("N116").Value = -5337

$ws.Range("H122").Value = 2321.64
$ws.Range("I122").Value = 2117.7727
$ws.Range("J122").Value = 3816.6667
$ws.Range("K122").Value = 6353.3181
$ws.Range("L122").Value = 11450.0001
$ws.Range("M122").Value = -3903.3181
$ws.Range("N122").Value = -16350.0001

$ws.Range("H132").Value = 4811.7856
$ws.Range("I132").Value = 2924.6875
$ws.Range("J132").Value = 7327.9165
$ws.Range("K132").Value = 8774.0625
$ws.Range("L132").Value = 21983.7495
$ws.Range("M132").Value = -6244.0625
$ws.Range("N132").Value = -27043.7495


# ===================== BSM =====================
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H3").Value = 1395.4286
$ws.Range("I3").Value = 1445.1538
$ws.Range("J3").Value = 749
$ws.Range("K3").Value = 1445.1538
$ws.Range("L3").Value = 749
$ws.Range("M3").Value = -1331.1538
$ws.Range("N3").Value = -977

$ws.Range("H20").Value = 1810.9584
$ws.Range("I20").Value = 1545
$ws.Range("K20").Value = 1545
$ws.Range("M20").Value = -1298

$ws.Range("H99").Value = 1936.875
$ws.Range("I99").Value = 1750
$ws.Range("J99").Value = 2497.5
$ws.Range("K99").Value = 1750
$ws.Range("L99").Value = 2497.5
$ws.Range("M99").Value = -252
$ws.Range("N99").Value = -5493.5

$ws.Range("H105").Value = 3668.6924
$ws.Range("I105").Value = 3502.9546
$ws.Range("J105").Value = 4580.25
$ws.Range("K105").Value = 3502.9546
$ws.Range("L105").Value = 4580.25
$ws.Range("M105").Value = -1755.9546
$ws.Range("N105").Value = -8074.25

$ws.Range("H107").Value = 3838.8064
$ws.Range("I107").Value = 3950.4
$ws.Range("J107").Value = 491
$ws.Range("K107").Value = 3950.4
$ws.Range("L107").Value = 491
$ws.Range("M107").Value = -2030.4
$ws.Range("N107").Value = -4331


# ===================== CRP =====================
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("I7").Value = 154.86363
$ws.Range("J7").Value = 506.66666
$ws.Range("K7").Value = 154.86363
$ws.Range("L7").Value = 506.66666
$ws.Range("M7").Value = -41.86363
$ws.Range("N7").Value = -732.66666

$ws.Range("H15").Value = 2227.5
$ws.Range("I15").Value = 1030
$ws.Range("J15").Value = 3082.8572
$ws.Range("K15").Value = 1030
$ws.Range("L15").Value = 3082.8572
$ws.Range("M15").Value = -860
$ws.Range("N15").Value = -3422.8572

$ws.Range("H16").Value = 719.63635
$ws.Range("I16").Value = 688.9375
$ws.Range("K16").Value = 688.9375
$ws.Range("M16").Value = -401.9375

$ws.Range("H31").Value = 5104.268
$ws.Range("I31").Value = 2002.7715
$ws.Range("K31").Value = 2002.7715
$ws.Range("M31").Value = -1707.7715

$ws.Range("H34").Value = 5104.268
$ws.Range("I34").Value = 2002.7715
$ws.Range("K34").Value = 2002.7715
$ws.Range("M34").Value = -1800.7715

$ws.Range("H39").Value = 37749
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 37749
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 37749
$ws.Range("M39").Value = $null
$ws.Range("N39").Value = -38531

$ws.Range("H49").Value = 37749
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 37749
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 37749
$ws.Range("M49").Value = $null
$ws.Range("N49").Value = -38113

$ws.Range("H60").Value = 24941.5
$ws.Range("I60").Value = 8412.5
$ws.Range("J60").Value = 57999.5
$ws.Range("K60").Value = 8412.5
$ws.Range("L60").Value = 57999.5
$ws.Range("M60").Value = -7901.5
$ws.Range("N60").Value = -59021.5

$ws.Range("H86").Value = 8131
$ws.Range("I86").Value = 7991.75
$ws.Range("J86").Value = 8316.666999999999
$ws.Range("K86").Value = 7991.75
$ws.Range("L86").Value = 8316.666999999999
$ws.Range("M86").Value = -6868.75
$ws.Range("N86").Value = -10562.667

$ws.Range("H89").Value = 8131
$ws.Range("I89").Value = 7991.75
$ws.Range("J89").Value = 8316.666999999999
$ws.Range("K89").Value = 39958.75
$ws.Range("L89").Value = 41583.335
$ws.Range("M89").Value = -34342.75
$ws.Range("N89").Value = -52815.335

$ws.Range("H105").Value = 2474.1667
$ws.Range("I105").Value = 2711.5
$ws.Range("J105").Value = 1999.5
$ws.Range("K105").Value = 2711.5
$ws.Range("L105").Value = 1999.5
$ws.Range("M105").Value = -964.5
$ws.Range("N105").Value = -5493.5

$ws.Range("H113").Value = 719.63635
$ws.Range("I113").Value = 688.9375
$ws.Range("K113").Value = 688.9375
$ws.Range("M113").Value = 1481.0625


# ===================== CUL =====================
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H5").Value = 1659
$ws.Range("I5").Value = 1737.6666
$ws.Range("J5").Value = 1501.6666
$ws.Range("K5").Value = 5212.9998
$ws.Range("L5").Value = 4504.9998
$ws.Range("M5").Value = -5100.9998
$ws.Range("N5").Value = -4728.9998

$ws.Range("H12").Value = 510.44446
$ws.Range("I12").Value = 150.5
$ws.Range("J12").Value = 613.2857
$ws.Range("K12").Value = 451.5
$ws.Range("L12").Value = 1839.8571
$ws.Range("M12").Value = -278.5
$ws.Range("N12").Value = -2185.8571

$ws.Range("H23").Value = 113
$ws.Range("J23").Value = 129.5
$ws.Range("L23").Value = 388.5
$ws.Range("N23").Value = -858.5

$ws.Range("H56").Value = 3643.2
$ws.Range("I56").Value = 3643.2
$ws.Range("K56").Value = 3643.2
$ws.Range("M56").Value = -3113.2

$ws.Range("H81").Value = 3496
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = $null

$ws.Range("H84").Value = 3496
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = $null

$ws.Range("H98").Value = 1366.4
$ws.Range("J98").Value = 1343.375
$ws.Range("L98").Value = 4030.125
$ws.Range("N98").Value = -7026.125

$ws.Range("H107").Value = 784.2857
$ws.Range("J107").Value = 784.2857
$ws.Range("L107").Value = 2352.8571
$ws.Range("N107").Value = -6192.8571

$ws.Range("H113").Value = 5435784
$ws.Range("J113").Value = 6212278.5
$ws.Range("L113").Value = 18636835.5
$ws.Range("N113").Value = -18641175.5

$ws.Range("H132").Value = 2936.4443
$ws.Range("J132").Value = 3586.25
$ws.Range("L132").Value = 32276.25
$ws.Range("N132").Value = -37336.25

$ws.Range("H133").Value = 4521.7812
$ws.Range("I133").Value = 3999.4
$ws.Range("J133").Value = 4618.5186
$ws.Range("K133").Value = 11998.2
$ws.Range("L133").Value = 13855.5558
$ws.Range("M133").Value = -6938.200000000001
$ws.Range("N133").Value = -23975.5558

$ws.Range("H135").Value = 1659
$ws.Range("I135").Value = 1737.6666
$ws.Range("J135").Value = 1501.6666
$ws.Range("K135").Value = 15638.9994
$ws.Range("L135").Value = 13514.9994
$ws.Range("M135").Value = -13103.9994
$ws.Range("N135").Value = -18584.9994

$ws.Range("H136").Value = 3985.9285
$ws.Range("I136").Value = 3441.9167
$ws.Range("K136").Value = 10325.7501
$ws.Range("M136").Value = -5225.750100000001

$ws.Range("H137").Value = 3676.64
$ws.Range("I137").Value = 2702.2222
$ws.Range("K137").Value = 8106.6666
$ws.Range("M137").Value = -3006.6666

$ws.Range("H140").Value = 1647.5319
$ws.Range("I140").Value = 1317.7391
$ws.Range("K140").Value = 3953.2173
$ws.Range("M140").Value = 1226.7827

$ws.Range("H141").Value = 3191.8333
$ws.Range("I141").Value = 3191.8333
$ws.Range("K141").Value = 9575.499899999999
$ws.Range("M141").Value = -4395.499899999999


# ===================== GSM =====================
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H43").Value = 29000
$ws.Range("I43").Value = 29000
$ws.Range("K43").Value = 29000
$ws.Range("M43").Value = -28849

$ws.Range("H49").Value = 5000
$ws.Range("I49").Value = 5000
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 5000
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -4816
$ws.Range("N49").Value = $null

$ws.Range("H58").Value = 15999.5
$ws.Range("I58").Value = 15999.5
$ws.Range("K58").Value = 15999.5
$ws.Range("M58").Value = -15722.5

$ws.Range("H80").Value = 5326
$ws.Range("I80").Value = 3625.2727
$ws.Range("J80").Value = 7998.5713
$ws.Range("K80").Value = 3625.2727
$ws.Range("L80").Value = 7998.5713
$ws.Range("M80").Value = -2627.2727
$ws.Range("N80").Value = -9994.5713

$ws.Range("H83").Value = 5326
$ws.Range("I83").Value = 3625.2727
$ws.Range("J83").Value = 7998.5713
$ws.Range("K83").Value = 18126.3635
$ws.Range("L83").Value = 39992.85649999999
$ws.Range("M83").Value = -13134.3635
$ws.Range("N83").Value = -49976.85649999999

$ws.Range("H102").Value = 2441.6667
$ws.Range("I102").Value = 1720.2222
$ws.Range("J102").Value = 5688.1665
$ws.Range("K102").Value = 1720.2222
$ws.Range("L102").Value = 5688.1665
$ws.Range("M102").Value = -98.22219999999993
$ws.Range("N102").Value = -8932.166499999999

$ws.Range("H107").Value = 1099.6666
$ws.Range("J107").Value = 1299
$ws.Range("L107").Value = 1299
$ws.Range("N107").Value = -5139

$ws.Range("H113").Value = 3761.2307
$ws.Range("I113").Value = 3717.9092
$ws.Range("K113").Value = 3717.9092
$ws.Range("M113").Value = -1547.9092

$ws.Range("H122").Value = 2538.2122
$ws.Range("I122").Value = 1311.4482
$ws.Range("K122").Value = 3934.3446
$ws.Range("M122").Value = -1484.3446

$ws.Range("H126").Value = 3720.3
$ws.Range("I126").Value = 3571.9565
$ws.Range("J126").Value = 4207.7144
$ws.Range("K126").Value = 10715.8695
$ws.Range("L126").Value = 12623.1432
$ws.Range("M126").Value = -8245.869499999999
$ws.Range("N126").Value = -17563.1432

$ws.Range("H133").Value = 83950
$ws.Range("J133").Value = 83950
$ws.Range("L133").Value = 83950
$ws.Range("N133").Value = -94070

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = $null


# ===================== LTW =====================
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H7").Value = 43885.6
$ws.Range("I7").Value = 45610.043
$ws.Range("K7").Value = 45610.043
$ws.Range("M7").Value = -45498.043

$ws.Range("H22").Value = 1740.1111
$ws.Range("I22").Value = 1718.75
$ws.Range("J22").Value = 1757.2
$ws.Range("K22").Value = 1718.75
$ws.Range("L22").Value = 1757.2
$ws.Range("M22").Value = -1423.75
$ws.Range("N22").Value = -2347.2

$ws.Range("H27").Value = 1740.1111
$ws.Range("I27").Value = 1718.75
$ws.Range("J27").Value = 1757.2
$ws.Range("K27").Value = 1718.75
$ws.Range("L27").Value = 1757.2
$ws.Range("M27").Value = -1611.75
$ws.Range("N27").Value = -1971.2

$ws.Range("H46").Value = 3277.8333
$ws.Range("I46").Value = 872.4
$ws.Range("K46").Value = 872.4
$ws.Range("M46").Value = -684.4

$ws.Range("H93").Value = 3060.8462
$ws.Range("I93").Value = 2161.5
$ws.Range("K93").Value = 2161.5
$ws.Range("M93").Value = -913.5

$ws.Range("H100").Value = 817.3333
$ws.Range("I100").Value = 817.3333
$ws.Range("K100").Value = 817.3333
$ws.Range("M100").Value = -276.3333

$ws.Range("H122").Value = 4010.9666
$ws.Range("I122").Value = 3516.2778
$ws.Range("K122").Value = 10548.8334
$ws.Range("M122").Value = -8098.8334

$ws.Range("H126").Value = 43885.6
$ws.Range("I126").Value = 45610.043
$ws.Range("K126").Value = 136830.129
$ws.Range("M126").Value = -134360.129

$ws.Range("H132").Value = 2212.9473
$ws.Range("I132").Value = 1401.931
$ws.Range("J132").Value = 4826.222
$ws.Range("K132").Value = 4205.793
$ws.Range("L132").Value = 14478.666
$ws.Range("M132").Value = -1675.793
$ws.Range("N132").Value = -19538.666

$ws.Range("H133").Value = 99974
$ws.Range("J133").Value = 99974
$ws.Range("L133").Value = 99974
$ws.Range("N133").Value = -105034

$ws.Range("H136").Value = 5685.2065
$ws.Range("I136").Value = 2848.8572
$ws.Range("K136").Value = 8546.571599999999
$ws.Range("M136").Value = -5996.571599999999

$ws.Range("H139").Value = 109499
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = $null

$ws.Range("H141").Value = 44570
$ws.Range("J141").Value = 44570
$ws.Range("L141").Value = 44570
$ws.Range("N141").Value = -54930


# ===================== WVR =====================
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H64").Value = 69967
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").Value = $null

$ws.Range("H67").Value = 69967
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").Value = $null

$ws.Range("H81").Value = 1317.0769
$ws.Range("I81").Value = 1317.0769
$ws.Range("K81").Value = 2634.1538
$ws.Range("M81").Value = -1573.1538

$ws.Range("H84").Value = 1317.0769
$ws.Range("I84").Value = 1317.0769
$ws.Range("K84").Value = 13170.769
$ws.Range("M84").Value = -7866.769

$ws.Range("H96").Value = 5102.25
$ws.Range("I96").Value = 3019.8
$ws.Range("J96").Value = 5796.4
$ws.Range("K96").Value = 3019.8
$ws.Range("L96").Value = 5796.4
$ws.Range("M96").Value = -1646.8
$ws.Range("N96").Value = -8542.4

$ws.Range("H132").Value = 1789.5763
$ws.Range("I132").Value = 1386.6296
$ws.Range("K132").Value = 4159.8888
$ws.Range("M132").Value = -1629.8888

$ws.Range("H136").Value = 1342.575
$ws.Range("I136").Value = 566.7879
$ws.Range("J136").Value = 4999.857
$ws.Range("K136").Value = 1700.3637
$ws.Range("L136").Value = 14999.571
$ws.Range("M136").Value = 849.6362999999999
$ws.Range("N136").Value = -20099.571

